$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Properties")

$ws.Range("C1").Value = "enabled"
$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $true
$ws.Range("C4").Value = $true

$ws.Columns.Item(2).ColumnWidth = 36.67

$ws.Range("C5").Select() | Out-Null
